$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# Update HU - 9 text (row 23)
$ws.Range("A23").Value = "HU - 9 Yo como usuario registrado, quiero finalizar el pedido/compra."

# Update HU - 10 text (row 24)
$ws.Range("A24").Value = "HU - 10 Yo como usuario registrado, puedo cambiar mi foto de perfil."

# Remove the old HU - 11 row entirely (row 25)
$ws.Rows("25").Delete()

# Update selection to match the new last populated cell
$ws.Range("A24").Select()
